$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.208092485549133
$ws.Range("C2").Value = 0.5028901734104047
$ws.Range("J2").Value = 0.008670520231213872
$ws.Range("P2").Value = 0.1820809248554913
$ws.Range("S2").Value = 0.09826589595375723
$ws.Range("C3").Value = 0.005405405405405406
$ws.Range("J3").Value = 0.05405405405405406
$ws.Range("P3").Value = 0.7567567567567568
$ws.Range("S3").Value = 0.1837837837837838
$ws.Range("J4").Value = 0.1052631578947368
$ws.Range("P4").Value = 0.6052631578947368
$ws.Range("S4").Value = 0.2894736842105263
$ws.Range("B6").Value = 0.01913875598086124
$ws.Range("D6").Value = 0.009569377990430622
$ws.Range("F6").Value = 0.06698564593301436
$ws.Range("J6").Value = 0.3301435406698565
$ws.Range("O6").Value = 0.02392344497607655
$ws.Range("Q6").Value = 0.1531100478468899
$ws.Range("R6").Value = 0.05263157894736842
$ws.Range("S6").Value = 0.3444976076555024
$ws.Range("B7").Value = 0.1534090909090909
$ws.Range("D7").Value = 0.005681818181818182
$ws.Range("F7").Value = 0.04545454545454546
$ws.Range("J7").Value = 0.1363636363636364
$ws.Range("O7").Value = 0.05113636363636364
$ws.Range("Q7").Value = 0.1590909090909091
$ws.Range("R7").Value = 0.05113636363636364
$ws.Range("S7").Value = 0.3977272727272727
$ws.Range("B8").Value = 0.09547738693467336
$ws.Range("D8").Value = 0.01758793969849246
$ws.Range("E8").Value = 0.002512562814070352
$ws.Range("F8").Value = 0.05276381909547739
$ws.Range("J8").Value = 0.1206030150753769
$ws.Range("O8").Value = 0.03266331658291458
$ws.Range("Q8").Value = 0.1658291457286432
$ws.Range("R8").Value = 0.1055276381909548
$ws.Range("S8").Value = 0.407035175879397
$ws.Range("B9").Value = 0.1147540983606557
$ws.Range("D9").Value = 0.0273224043715847
$ws.Range("E9").Value = 0.00546448087431694
$ws.Range("F9").Value = 0.04918032786885246
$ws.Range("J9").Value = 0.1202185792349727
$ws.Range("O9").Value = 0.01092896174863388
$ws.Range("Q9").Value = 0.1311475409836066
$ws.Range("R9").Value = 0.1147540983606557
$ws.Range("S9").Value = 0.4262295081967213
$ws.Range("B10").Value = 0.125948406676783
$ws.Range("D10").Value = 0.01820940819423369
$ws.Range("E10").Value = 0.0007587253414264037
$ws.Range("F10").Value = 0.0629742033383915
$ws.Range("J10").Value = 0.1350531107738998
$ws.Range("O10").Value = 0.02579666160849772
$ws.Range("Q10").Value = 0.1965098634294385
$ws.Range("R10").Value = 0.07814871016691957
$ws.Range("S10").Value = 0.3566009104704097
$ws.Range("G11").Value = 0.1863799283154122
$ws.Range("J11").Value = 0.07885304659498207
$ws.Range("K11").Value = 0.2293906810035842
$ws.Range("L11").Value = 0.4982078853046595
$ws.Range("S11").Value = 0.007168458781362007
$ws.Range("G12").Value = 0.7272727272727273
$ws.Range("J12").Value = 0.2237762237762238
$ws.Range("K12").Value = 0.01398601398601399
$ws.Range("L12").Value = 0.01398601398601399
$ws.Range("S12").Value = 0.02097902097902098
$ws.Range("G13").Value = 0.6666666666666666
$ws.Range("J13").Value = 0.2142857142857143
$ws.Range("S13").Value = 0.119047619047619
$ws.Range("F15").Value = 0.04347826086956522
$ws.Range("H15").Value = 0.1521739130434783
$ws.Range("I15").Value = 0.06086956521739131
$ws.Range("J15").Value = 0.3695652173913043
$ws.Range("K15").Value = 0.0391304347826087
$ws.Range("M15").Value = 0.004347826086956522
$ws.Range("O15").Value = 0.05217391304347826
$ws.Range("S15").Value = 0.2782608695652174
$ws.Range("F16").Value = 0.02702702702702703
$ws.Range("H16").Value = 0.1351351351351351
$ws.Range("I16").Value = 0.07207207207207207
$ws.Range("J16").Value = 0.3828828828828829
$ws.Range("K16").Value = 0.1126126126126126
$ws.Range("M16").Value = 0.03153153153153153
$ws.Range("O16").Value = 0.04054054054054054
$ws.Range("S16").Value = 0.1981981981981982
$ws.Range("F17").Value = 0.02891566265060241
$ws.Range("H17").Value = 0.1493975903614458
$ws.Range("I17").Value = 0.0963855421686747
$ws.Range("J17").Value = 0.4626506024096386
$ws.Range("K17").Value = 0.0891566265060241
$ws.Range("M17").Value = 0.01686746987951807
$ws.Range("O17").Value = 0.06265060240963856
$ws.Range("S17").Value = 0.09397590361445783
$ws.Range("F18").Value = 0.02702702702702703
$ws.Range("H18").Value = 0.1945945945945946
$ws.Range("I18").Value = 0.1351351351351351
$ws.Range("J18").Value = 0.3945945945945946
$ws.Range("K18").Value = 0.04864864864864865
$ws.Range("M18").Value = 0.01621621621621622
$ws.Range("O18").Value = 0.03783783783783784
$ws.Range("S18").Value = 0.145945945945946
$ws.Range("F19").Value = 0.02093397745571659
$ws.Range("H19").Value = 0.1948470209339775
$ws.Range("I19").Value = 0.06763285024154589
$ws.Range("J19").Value = 0.3985507246376812
$ws.Range("K19").Value = 0.106280193236715
$ws.Range("M19").Value = 0.02012882447665056
$ws.Range("O19").Value = 0.06763285024154589
$ws.Range("S19").Value = 0.1239935587761675
